# Add new translation rows to the "Import" sheet (sheet1), mirroring the
# existing layout/style of the last existing row (row 344).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Import")

$rows = @(
    @("cs", "lab.base.tooltip.create", "Nová báze"),
    @("cs", "lab.base.tooltip.create", "Nová báze"),
    @("cs", "lab.base.create.title", "Nová báze"),
    @("cs", "lab.base.create.subtitle", "Báze je užitečná v mixech."),
    @("cs", "lab.base.name.label", "Název báze"),
    @("cs", "lab.base.vendorId.label", "Výrobce"),
    @("cs", "lab.base.pg.label", "PG"),
    @("cs", "lab.base.vg.label", "VG"),
    @("cs", "lab.base.create.submit", "Vytvořit bázi"),
    @("cs", "lab.base.create.success", "Báze [{{data.name}}] byla uložena.")
)

$startRow = 345
$lastExistingRow = 344

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i

    # Copy formatting from the last existing data row so the new row keeps
    # the same style (wrap text etc.) without creating duplicate style
    # entries.
    $ws.Range("A$lastExistingRow`:C$lastExistingRow").Copy()
    $ws.Range("A$r`:C$r").PasteSpecial(-4122)

    $data = $rows[$i]
    $ws.Range("A$r").Value = $data[0]
    $ws.Range("B$r").Value = $data[1]
    $ws.Range("C$r").Value = $data[2]
}

$excel.CutCopyMode = 0

# Update the sheet view selection to match the new last entry.
[void]$ws.Range("B353").Select()

Write-Output "Added $($rows.Length) rows"
